# Add newly-gathered headshot rows (29-38) to the "headshot_url" tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$logoUrl = "https://images.squarespace-cdn.com/content/v1/5f63780dbc8d16716cca706a/1604523297465-6BAIW9AOVGRI7PARBCH3/rowing-canada-new-logo.jpg"
$lauraCourtUrl = "https://rowingcanada.org/uploads/2018/11/Laura-Court-Rowing-Canada-Sep-1-2022-%C2%A9KevinLightPhoto-_RL_8664A-1920x1869.jpg"

# --- Row 29: Grace Vanden Broek ---
$ws.Range("A29").Value = "Grace Vanden Broek"
$ws.Range("B29").Value = $logoUrl

# --- Row 30: Axel Ewashko ---
$ws.Range("A30").Value = "Axel Ewashko"
$ws.Range("B30").Value = $logoUrl

# --- Row 31: Nicole Cusack ---
$ws.Range("A31").Value = "Nicole Cusack"
$ws.Range("B31").Value = $logoUrl

# --- Row 32: Alizee Brien ---
$ws.Range("A32").Value = "Alizee Brien"
$ws.Range("B32").Value = $logoUrl

# --- Row 33: Mitchell Rodgers ---
$ws.Range("A33").Value = "Mitchell Rodgers"
$ws.Range("B33").Value = $logoUrl

# --- Row 34: Brenna Randall ---
$ws.Range("A34").Value = "Brenna Randall"
$ws.Range("B34").Value = $logoUrl

# --- Row 35: Laura Court (link entered before the name here) ---
$ws.Range("B35").Value = $lauraCourtUrl
$ws.Range("A35").Value = "Laura Court"

# --- Row 36: Trevor Jones ---
$ws.Range("A36").Value = "Trevor Jones"
$ws.Range("B36").Value = $logoUrl

# --- Row 37: Alexis Cronk ---
$ws.Range("A37").Value = "Alexis Cronk"
$ws.Range("B37").Value = $logoUrl

# --- Row 38: Piper Battersby ---
$ws.Range("A38").Value = "Piper Battersby"
$ws.Range("B38").Value = $logoUrl

# Hyperlink the "no confirmed photo yet" / found-link cells, in the order
# they were actually wired up (B35 first, then B34, then B37, then B38).
$ws.Hyperlinks.Add($ws.Range("B35"), $lauraCourtUrl) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B34"), $logoUrl) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B37"), $logoUrl) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B38"), $logoUrl) | Out-Null

$ws.Range("D42").Select()
